$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prepare the LAST SCRAPE DATE column (F) so the date-like strings are
# written as plain text instead of being auto-converted into date serials.
$dateRange = $ws.Range("F2:F51")
$dateRange.NumberFormat = "@"

$data = @(
    @('WV', '$1.00 Games', 'Bonus Cash', '957', '1', '2019-03-12'),
    @('WV', '$1.00 Games', 'Veterans'' Cash', '990', '2', '2019-03-12'),
    @('WV', '$1.00 Games', 'Quick 6', '967', '4', '2019-03-12'),
    @('WV', '$1.00 Games', 'Tic Tac Toe', '950', '1', '2019-03-12'),
    @('WV', '$1.00 Games', 'Black', '999', '3', '2019-02-25'),
    @('WV', '$1.00 Games', 'SCARDEY CASH / SANTA PAWS', '978', '1', '2019-03-12'),
    @('WV', '$1.00 Games', 'Big Cheese', '1003', '7', '2019-02-22'),
    @('WV', '$1.00 Games', '3X', '974', '0', '2019-03-12'),
    @('WV', '$2.00 Games', 'High 5s', '958', '1', '2019-03-12'),
    @('WV', '$2.00 Games', 'LOS MUERTOS / CASHING SNOW', '977', '4', '2019-03-12'),
    @('WV', '$2.00 Games', 'Gem 7s', '993', '5', '2019-03-12'),
    @('WV', '$2.00 Games', '10 Grand', '1032', '2', '2019-03-12'),
    @('WV', '$2.00 Games', 'Bacon Love - Jack Cash', '1002', '4', '2019-03-07'),
    @('WV', '$2.00 Games', 'Keno Bullseye', '970', '1', '2019-03-12'),
    @('WV', '$2.00 Games', 'Classic Black', '998', '2', '2019-03-07'),
    @('WV', '$2.00 Games', 'Winner Winner Chicken Dinner', '961', '161', '2019-03-12'),
    @('WV', '$2.00 Games', '5X', '973', '0', '2019-03-12'),
    @('WV', '$2.00 Games', '7s', '940', '0', '2019-03-12'),
    @('WV', '$2.00 Games', 'Red White & Blue', '962', '0', '2019-03-12'),
    @('WV', '$3.00 Games', 'Cash Wheel', '963', '1', '2019-03-12'),
    @('WV', '$3.00 Games', 'Clockwork Crossword', '959', '1', '2019-03-12'),
    @('WV', '$3.00 Games', '8 Ball Bingo', '989', '3', '2019-03-12'),
    @('WV', '$3.00 Games', 'CASH CONNECT 1007', '1007', '3', '2019-03-12'),
    @('WV', '$3.00 Games', 'ROLLIN'' BIG', '992', '1', '2019-03-12'),
    @('WV', '$3.00 Games', 'Hot Chili Crossword', '979', '1', '2019-03-12'),
    @('WV', '$3.00 Games', 'Winning Numbers', '966', '1', '2019-02-11'),
    @('WV', '$3.00 Games', 'Blingo', '935', '0', '2019-03-12'),
    @('WV', '$5.00 Games', 'Red & Green Machine', '976', '1', '2019-03-12'),
    @('WV', '$5.00 Games', 'Crack the Code', '942', '1', '2019-03-12'),
    @('WV', '$5.00 Games', 'Cash Squatch', '1000', '2', '2019-03-12'),
    @('WV', '$5.00 Games', '10X', '972', '2', '2019-03-12'),
    @('WV', '$5.00 Games', 'Mega Black', '997', '3', '2019-03-12'),
    @('WV', '$5.00 Games', 'VIVA LAS KENO 1005', '1005', '3', '2019-03-12'),
    @('WV', '$5.00 Games', 'Jumbo Jack', '994', '2', '2019-03-12'),
    @('WV', '$5.00 Games', 'PAC-MAN® & Ms. PAC-MAN®', '964', '1', '2019-03-12'),
    @('WV', '$5.00 Games', 'Loaded', '947', '14', '2019-03-12'),
    @('WV', '$5.00 Games', '$40 GRAND', '968', '2', '2019-03-12'),
    @('WV', '$5.00 Games', 'More Money', '960', '1', '2019-03-12'),
    @('WV', '$5.00 Games', 'Cash Extravaganza', '931', '0', '2019-03-12'),
    @('WV', '$5.00 Games', 'Cherries Wild', '955', '0', '2019-03-12'),
    @('WV', '$5.00 Games', 'Viva Las Keno', '969', '0', '2019-03-12'),
    @('WV', '$5.00 Games', 'TRICK OR TREAT', '975', '0', '2019-03-12'),
    @('WV', '$10.00 Games', 'Flawless Fortune', '965', '2', '2019-03-12'),
    @('WV', '$10.00 Games', 'Quick 10 Bonus Jackpot', '956', '5', '2019-03-12'),
    @('WV', '$10.00 Games', 'Maximum Black', '991', '3', '2019-03-12'),
    @('WV', '$10.00 Games', '$50 or $100', '995', '4384', '2019-03-12'),
    @('WV', '$10.00 Games', '$100 $200 or $500', '1001', '72', '2019-03-12'),
    @('WV', '$10.00 Games', 'Extreme Black', '996', '2', '2019-03-12'),
    @('WV', '$10.00 Games', '20X', '971', '0', '2019-03-12'),
    @('WV', '$10.00 Games', '20X Cash / $5,000 Blowout', '949', '0', '2019-03-12')
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = 2 + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = [double]$row[3]
    $ws.Cells.Item($r, 5).Value = [double]$row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

# Remove the temporary text format now that the values are stored so the
# cells fall back to the default/general style, matching the original file.
$dateRange.ClearFormats()

Write-Host "Done"
